# feat: add 2022-Q1 data
#
# The workbook currently ends with a "总计" (totals) summary sheet.
# This change:
#   1. Turns that existing "总计" sheet into the new "2022-Q1" holdings
#      sheet (same layout as the other quarterly sheets: 基金代码 / 基金名称 /
#      基金规模 / 股票总仓位 / 仓位占比 / 持有市值(亿元) / 仓位排名).
#   2. Appends a brand-new "总计" sheet after it, containing the same
#      summary table as before plus a new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# Style source: a cell on an existing quarterly sheet that already uses
# the bold/centered/bordered header-and-index style (internal style index 2).
$styleSrc = $wb.Worksheets.Item("2021-Q4").Range("B1")

# A cell with the plain/default style (index 0) - used to strip away any
# number-format-only style we create below when forcing text-typed values.
$plainSrc = $wb.Worksheets.Item("2021-Q4").Range("Z100")

# ---------------------------------------------------------------------
# Step 1: insert the new "总计" sheet right after the current one, then
# rename sheets so the old "总计" becomes "2022-Q1" and the new sheet
# becomes "总计" (matches the workbook.xml sheetId/order in the diff).
# ---------------------------------------------------------------------
$oldTotal = $wb.Worksheets.Item("总计")
$newTotal = $wb.Worksheets.Add($null, $oldTotal)

$newTotal.PageSetup.LeftMargin = 54
$newTotal.PageSetup.RightMargin = 54
$newTotal.PageSetup.TopMargin = 72
$newTotal.PageSetup.BottomMargin = 72
$newTotal.PageSetup.HeaderMargin = 36
$newTotal.PageSetup.FooterMargin = 36
$newTotal.Outline.SummaryRow = 1
$newTotal.Outline.SummaryColumn = 1

$oldTotal.Name = "2022-Q1"
$newTotal.Name = "总计"

$q1 = $wb.Worksheets.Item("2022-Q1")
$tot = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# Step 2: rebuild "2022-Q1" with the fund-holdings table (it currently
# still holds the old totals-table content, so clear it first).
# ---------------------------------------------------------------------
$q1.Cells.Clear()

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"
$styleSrc.Copy()
$q1.Range("B1:H1").PasteSpecial($xlPasteFormats)

$q1Rows = @(
    @(0, "501021", "华宝兴业标普香港上市中国中小盘指数(QDII-LOF)A", "4.96", "94.77", "1.62", "0.0804", 9),
    @(1, "004532", "民生加银中证港股通高股息精选指数A",             "0.26", "94.88", "3.57", "0.0093", 8),
    @(2, "006127", "华宝兴业标普香港上市中国中小盘指数(QDII-LOF)C", "0.23", "94.77", "1.62", "0.0037", 9),
    @(3, "011647", "博时港股通红利精选混合A",                       "0.13", "92.10", "2.84", "0.0037", 10),
    @(4, "004533", "民生加银中证港股通高股息精选指数C",             "0.10", "94.88", "3.57", "0.0036", 8),
    @(5, "011648", "博时港股通红利精选混合C",                       "0.02", "92.10", "2.84", "0.0006", 10),
    @(6, "005770", "信达澳银中证沪港深高股息精选指数",               "0.01", "92.47", "2.63", "0.0003", 4)
)

$r = 2
foreach ($row in $q1Rows) {
    $q1.Cells.Item($r, 1).Value = $row[0]

    # Fund code looks numeric but must stay text (keeps leading zeros).
    $q1.Cells.Item($r, 2).NumberFormat = "@"
    $q1.Cells.Item($r, 2).Value = $row[1]

    $q1.Cells.Item($r, 3).Value = $row[2]

    $q1.Cells.Item($r, 4).NumberFormat = "@"
    $q1.Cells.Item($r, 4).Value = $row[3]
    $q1.Cells.Item($r, 5).NumberFormat = "@"
    $q1.Cells.Item($r, 5).Value = $row[4]
    $q1.Cells.Item($r, 6).NumberFormat = "@"
    $q1.Cells.Item($r, 6).Value = $row[5]
    $q1.Cells.Item($r, 7).NumberFormat = "@"
    $q1.Cells.Item($r, 7).Value = $row[6]

    $q1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Restore the default (no explicit number-format) style on the text-forced
# B and D:G cells while keeping their stored value typed as text.
$plainSrc.Copy()
$q1.Range("B2:B8").PasteSpecial($xlPasteFormats)
$plainSrc.Copy()
$q1.Range("D2:G8").PasteSpecial($xlPasteFormats)

# Re-apply the bold/centered/bordered style to column A (index numbers).
$styleSrc.Copy()
$q1.Range("A2:A8").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Step 3: rebuild "总计" with the quarterly totals table, adding the new
# 2022-Q1 row at the top (index 0) and shifting the rest down.
# ---------------------------------------------------------------------
$tot.Range("B1").Value = "日期"
$tot.Range("C1").Value = "持有数量(只)"
$tot.Range("D1").Value = "持有市值(亿元)"
$styleSrc.Copy()
$tot.Range("B1:D1").PasteSpecial($xlPasteFormats)

$totRows = @(
    @(0, "2022-Q1", 7, 0.1),
    @(1, "2021-Q4", 5, 0.36),
    @(2, "2021-Q3", 17, 3.04),
    @(3, "2021-Q2", 33, 17.24),
    @(4, "2021-Q1", 25, 11.75),
    @(5, "2020-Q4", 24, 3.01)
)

$r = 2
foreach ($row in $totRows) {
    $tot.Cells.Item($r, 1).Value = $row[0]
    $tot.Cells.Item($r, 2).Value = $row[1]
    $tot.Cells.Item($r, 3).Value = $row[2]
    $tot.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$styleSrc.Copy()
$tot.Range("A2:A7").PasteSpecial($xlPasteFormats)

Write-Output ($wb.Worksheets | ForEach-Object { $_.Name } | Out-String)
